$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# New column A = BrowserName; existing "Categories" data moves to column B; new column C = URL
$ws.Range("A1").Value = "BrowserName"
$ws.Range("B1").Value = "Categories"
$ws.Range("C1").Value = "URL"

$ws.Range("A2").Value = "CHROME"
$ws.Range("B2").Value = "Laptops"
$ws.Range("C2").Value = "https://sapui5.hana.ondemand.com/test-resources/sap/m/demokit/cart/webapp/index.html"

$ws.Range("A3").Value = "CHROME"
$ws.Range("B3").Value = "Accessories"
$ws.Range("C3").Value = "https://sapui5.hana.ondemand.com/test-resources/sap/m/demokit/cart/webapp/index.html"

$ws.Range("C2:C3").Style = $ws.Range("A2").Style

$ws.Columns.Item(1).ColumnWidth = 12.859375
$ws.Columns.Item(3).ColumnWidth = 76.16796875

$ws.Range("C4").Select()
